$wb = $excel.ActiveWorkbook

# ---- Sheet1 (总计): insert new row for 2022-Q3, shift existing quarters down ----
$summary = $wb.Worksheets.Item(1)
$summary.Rows.Item(2).Insert()
$summary.Range("A2:D2").ClearFormats()
$summary.Range("B1").Copy()
$summary.Range("A2").PasteSpecial(-4122)

$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3
$summary.Range("A6").Value = 4
$summary.Range("A7").Value = 5

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 25
$summary.Range("D2").Value = 8.039999999999999

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 9
$summary.Range("D3").Value = 2.78

$summary.Range("B4").Value = "2022-Q1"
$summary.Range("C4").Value = 8
$summary.Range("D4").Value = 2.41

$summary.Range("B5").Value = "2021-Q4"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.42

$summary.Range("B6").Value = "2021-Q3"
$summary.Range("C6").Value = 4
$summary.Range("D6").Value = 0.11

$summary.Range("B7").Value = "2021-Q2"
$summary.Range("C7").Value = 2
$summary.Range("D7").Value = 0

# ---- New sheet "2022-Q3": fund holdings detail, inserted right after 总计 ----
$afterSheet = $wb.Worksheets.Item(2)
$newSheet = $wb.Worksheets.Add($afterSheet)
$newSheet.Name = "2022-Q3"
$ws = $newSheet

$ws.Range("B1").Value = "基金代码"
$ws.Range("C1").Value = "基金名称"
$ws.Range("D1").Value = "基金规模"
$ws.Range("E1").Value = "股票总仓位"
$ws.Range("F1").Value = "仓位占比"
$ws.Range("G1").Value = "持有市值(亿元)"
$ws.Range("H1").Value = "仓位排名"

$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'006113"
$ws.Range("C2").Value = "汇添富创新医药主题混合"
$ws.Range("D2").Value = "'85.14"
$ws.Range("E2").Value = "'75.99"
$ws.Range("F2").Value = "'2.92"
$ws.Range("G2").Value = "'2.4861"
$ws.Range("H2").Value = 10

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'001417"
$ws.Range("C3").Value = "汇添富医疗服务灵活配置混合A"
$ws.Range("D3").Value = "'32.06"
$ws.Range("E3").Value = "'85.81"
$ws.Range("F3").Value = "'4.22"
$ws.Range("G3").Value = "'1.3529"
$ws.Range("H3").Value = 5

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'000452"
$ws.Range("C4").Value = "南方医药保健灵活配置混合A"
$ws.Range("D4").Value = "'28.31"
$ws.Range("E4").Value = "'93.36"
$ws.Range("F4").Value = "'4.61"
$ws.Range("G4").Value = "'1.3051"
$ws.Range("H4").Value = 7

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'470006"
$ws.Range("C5").Value = "汇添富医药保健混合A"
$ws.Range("D5").Value = "'46.47"
$ws.Range("E5").Value = "'84.46"
$ws.Range("F5").Value = "'2.78"
$ws.Range("G5").Value = "'1.2919"
$ws.Range("H5").Value = 10

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'009664"
$ws.Range("C6").Value = "汇添富医疗积极成长一年持有期混合A"
$ws.Range("D6").Value = "'29.13"
$ws.Range("E6").Value = "'66.68"
$ws.Range("F6").Value = "'3.13"
$ws.Range("G6").Value = "'0.9118"
$ws.Range("H6").Value = 9

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'010054"
$ws.Range("C7").Value = "万家健康产业混合A"
$ws.Range("D7").Value = "'4.60"
$ws.Range("E7").Value = "'86.96"
$ws.Range("F7").Value = "'3.30"
$ws.Range("G7").Value = "'0.1518"
$ws.Range("H7").Value = 10

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'630010"
$ws.Range("C8").Value = "华商价值精选混合"
$ws.Range("D8").Value = "'4.30"
$ws.Range("E8").Value = "'81.81"
$ws.Range("F8").Value = "'2.94"
$ws.Range("G8").Value = "'0.1264"
$ws.Range("H8").Value = 10

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'009665"
$ws.Range("C9").Value = "汇添富医疗积极成长一年持有期混合C"
$ws.Range("D9").Value = "'3.63"
$ws.Range("E9").Value = "'66.68"
$ws.Range("F9").Value = "'3.13"
$ws.Range("G9").Value = "'0.1136"
$ws.Range("H9").Value = 9

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'506009"
$ws.Range("C10").Value = "国泰科创板两年定期开放混合"
$ws.Range("D10").Value = "'2.05"
$ws.Range("E10").Value = "'85.80"
$ws.Range("F10").Value = "'3.73"
$ws.Range("G10").Value = "'0.0765"
$ws.Range("H10").Value = 7

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'008107"
$ws.Range("C11").Value = "华商医药医疗行业股票"
$ws.Range("D11").Value = "'1.39"
$ws.Range("E11").Value = "'89.19"
$ws.Range("F11").Value = "'4.05"
$ws.Range("G11").Value = "'0.0563"
$ws.Range("H11").Value = 8

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'012358"
$ws.Range("C12").Value = "汇丰晋信医疗先锋混合A"
$ws.Range("D12").Value = "'1.87"
$ws.Range("E12").Value = "'58.14"
$ws.Range("F12").Value = "'2.52"
$ws.Range("G12").Value = "'0.0471"
$ws.Range("H12").Value = 7

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'620001"
$ws.Range("C13").Value = "金元顺安宝石动力混合"
$ws.Range("D13").Value = "'1.01"
$ws.Range("E13").Value = "'40.12"
$ws.Range("F13").Value = "'2.67"
$ws.Range("G13").Value = "'0.0270"
$ws.Range("H13").Value = 7

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "'010055"
$ws.Range("C14").Value = "万家健康产业混合C"
$ws.Range("D14").Value = "'0.78"
$ws.Range("E14").Value = "'86.96"
$ws.Range("F14").Value = "'3.30"
$ws.Range("G14").Value = "'0.0257"
$ws.Range("H14").Value = 10

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "'630006"
$ws.Range("C15").Value = "华商产业升级混合"
$ws.Range("D15").Value = "'0.85"
$ws.Range("E15").Value = "'81.97"
$ws.Range("F15").Value = "'2.94"
$ws.Range("G15").Value = "'0.0250"
$ws.Range("H15").Value = 10

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "'005493"
$ws.Range("C16").Value = "鑫元价值精选灵活配置混合A"
$ws.Range("D16").Value = "'0.55"
$ws.Range("E16").Value = "'76.82"
$ws.Range("F16").Value = "'3.36"
$ws.Range("G16").Value = "'0.0185"
$ws.Range("H16").Value = 2

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "'014933"
$ws.Range("C17").Value = "南方医药保健灵活配置混合C"
$ws.Range("D17").Value = "'0.13"
$ws.Range("E17").Value = "'93.36"
$ws.Range("F17").Value = "'4.61"
$ws.Range("G17").Value = "'0.0060"
$ws.Range("H17").Value = 7

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "'014285"
$ws.Range("C18").Value = "鑫元健康产业混合A"
$ws.Range("D18").Value = "'0.12"
$ws.Range("E18").Value = "'78.73"
$ws.Range("F18").Value = "'3.61"
$ws.Range("G18").Value = "'0.0043"
$ws.Range("H18").Value = 3

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "'006193"
$ws.Range("C19").Value = "鑫元核心资产股票A"
$ws.Range("D19").Value = "'0.11"
$ws.Range("E19").Value = "'83.48"
$ws.Range("F19").Value = "'3.18"
$ws.Range("G19").Value = "'0.0035"
$ws.Range("H19").Value = 7

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "'012359"
$ws.Range("C20").Value = "汇丰晋信医疗先锋混合C"
$ws.Range("D20").Value = "'0.12"
$ws.Range("E20").Value = "'58.14"
$ws.Range("F20").Value = "'2.52"
$ws.Range("G20").Value = "'0.0030"
$ws.Range("H20").Value = 7

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "'014286"
$ws.Range("C21").Value = "鑫元健康产业混合C"
$ws.Range("D21").Value = "'0.07"
$ws.Range("E21").Value = "'78.73"
$ws.Range("F21").Value = "'3.61"
$ws.Range("G21").Value = "'0.0025"
$ws.Range("H21").Value = 3

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "'015121"
$ws.Range("C22").Value = "汇添富医疗服务灵活配置混合C"
$ws.Range("D22").Value = "'0.02"
$ws.Range("E22").Value = "'85.81"
$ws.Range("F22").Value = "'4.22"
$ws.Range("G22").Value = "'0.0008"
$ws.Range("H22").Value = 5

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "'005494"
$ws.Range("C23").Value = "鑫元价值精选灵活配置混合C"
$ws.Range("D23").Value = "'0.01"
$ws.Range("E23").Value = "'76.82"
$ws.Range("F23").Value = "'3.36"
$ws.Range("G23").Value = "'0.0003"
$ws.Range("H23").Value = 2

$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "'006194"
$ws.Range("C24").Value = "鑫元核心资产股票C"
$ws.Range("D24").Value = "'0.01"
$ws.Range("E24").Value = "'83.48"
$ws.Range("F24").Value = "'3.18"
$ws.Range("G24").Value = "'0.0003"
$ws.Range("H24").Value = 7

$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "'960015"
$ws.Range("C25").Value = "汇添富医药保健混合O"
$ws.Range("D25").Value = "'0.00"
$ws.Range("E25").Value = "'84.46"
$ws.Range("F25").Value = "'2.78"
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 10

$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "'015122"
$ws.Range("C26").Value = "汇添富医疗服务灵活配置混合D"
$ws.Range("D26").Value = "'0.00"
$ws.Range("E26").Value = "'85.81"
$ws.Range("F26").Value = "'4.22"
$ws.Range("G26").Value = 0
$ws.Range("H26").Value = 5

# Reset auto-applied quote-prefix style on numeric-looking text cells
$ws.Range("B2").Style = "Normal"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Style = "Normal"
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Style = "Normal"
$ws.Range("B3").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Style = "Normal"
$ws.Range("F3").Style = "Normal"
$ws.Range("G3").Style = "Normal"
$ws.Range("B4").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Style = "Normal"
$ws.Range("F4").Style = "Normal"
$ws.Range("G4").Style = "Normal"
$ws.Range("B5").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Style = "Normal"
$ws.Range("F5").Style = "Normal"
$ws.Range("G5").Style = "Normal"
$ws.Range("B6").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Style = "Normal"
$ws.Range("G6").Style = "Normal"
$ws.Range("B7").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Style = "Normal"
$ws.Range("G7").Style = "Normal"
$ws.Range("B8").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Style = "Normal"
$ws.Range("F8").Style = "Normal"
$ws.Range("G8").Style = "Normal"
$ws.Range("B9").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Style = "Normal"
$ws.Range("F9").Style = "Normal"
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Style = "Normal"
$ws.Range("F10").Style = "Normal"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Style = "Normal"
$ws.Range("F11").Style = "Normal"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Style = "Normal"
$ws.Range("F12").Style = "Normal"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Style = "Normal"
$ws.Range("F13").Style = "Normal"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Style = "Normal"
$ws.Range("F14").Style = "Normal"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Style = "Normal"
$ws.Range("F15").Style = "Normal"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Style = "Normal"
$ws.Range("F16").Style = "Normal"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Style = "Normal"
$ws.Range("F17").Style = "Normal"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Style = "Normal"
$ws.Range("F18").Style = "Normal"
$ws.Range("G18").Style = "Normal"
$ws.Range("B19").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Style = "Normal"
$ws.Range("F19").Style = "Normal"
$ws.Range("G19").Style = "Normal"
$ws.Range("B20").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Style = "Normal"
$ws.Range("F20").Style = "Normal"
$ws.Range("G20").Style = "Normal"
$ws.Range("B21").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Style = "Normal"
$ws.Range("F21").Style = "Normal"
$ws.Range("G21").Style = "Normal"
$ws.Range("B22").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Style = "Normal"
$ws.Range("F22").Style = "Normal"
$ws.Range("G22").Style = "Normal"
$ws.Range("B23").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Style = "Normal"
$ws.Range("F23").Style = "Normal"
$ws.Range("G23").Style = "Normal"
$ws.Range("B24").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Style = "Normal"
$ws.Range("F24").Style = "Normal"
$ws.Range("G24").Style = "Normal"
$ws.Range("B25").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Style = "Normal"
$ws.Range("F25").Style = "Normal"
$ws.Range("B26").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Style = "Normal"
$ws.Range("F26").Style = "Normal"

# Apply header (bold+border) style and index-column style from a sibling quarter sheet
$template = $wb.Worksheets.Item(3)
$template.Range("B1:H1").Copy()
$ws.Range("B1:H1").PasteSpecial(-4122)
$template.Range("A2").Copy()
$ws.Range("A2:A26").PasteSpecial(-4122)

